$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update parameter values
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 48

# Update the active cell selection
$ws.Range("B3").Select()
